$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $styleDonor, $value) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
    $ws.Range($addr).Style = $ws.Range($styleDonor).Style
}

Set-TextCell "D2" "B2" '66.642.33'
$ws.Range("E2").Value = '  +0.29%  '

Set-TextCell "D3" "B3" '3.228.09'
$ws.Range("E3").Value = '  +0.97%  '

Set-TextCell "D4" "B4" '1.00'
$ws.Range("E4").Value = '  -0.08%  '

Set-TextCell "D5" "B5" '608.87'
$ws.Range("E5").Value = '  +2.09%  '

Set-TextCell "D6" "B6" '158.93'
$ws.Range("E6").Value = '  +2.86%  '

Set-TextCell "D7" "B7" '1.00'
$ws.Range("E7").Value = '  +0.00%  '

Set-TextCell "D8" "B8" '3.229.74'
$ws.Range("E8").Value = '  +1.03%  '

Set-TextCell "D9" "B9" '0.551'
$ws.Range("E9").Value = '  +1.04%  '

Set-TextCell "D10" "B10" '0.161'
$ws.Range("E10").Value = '  +0.23%  '

Set-TextCell "D11" "B11" '5.72'
$ws.Range("E11").Value = '  -5.15%  '

Set-TextCell "D12" "B12" '0.503'
$ws.Range("E12").Value = '  -2.93%  '

$ws.Range("E13").Value = '  +1.12%  '

Set-TextCell "D14" "B14" '38.81'
$ws.Range("E14").Value = '  -0.93%  '

Set-TextCell "D15" "B15" '3.760.16'
$ws.Range("E15").Value = '  +0.98%  '

Set-TextCell "D16" "B16" '66.692.94'
$ws.Range("E16").Value = '  +0.35%  '

Set-TextCell "D17" "B17" '7.36'
$ws.Range("E17").Value = '  -1.40%  '

Set-TextCell "D18" "B18" '3.235.93'
$ws.Range("E18").Value = '  +1.20%  '

$ws.Range("E19").Value = '  +1.08%  '

Set-TextCell "D20" "B20" '507.08'
$ws.Range("E20").Value = '  -1.60%  '

Set-TextCell "D21" "B21" '15.18'
$ws.Range("E21").Value = '  -1.54%  '

Set-TextCell "D22" "B22" '0.734'
$ws.Range("E22").Value = '  -1.03%  '

Set-TextCell "D23" "B23" '8.00'
$ws.Range("E23").Value = '  -1.66%  '

Set-TextCell "D24" "B24" '14.57'
$ws.Range("E24").Value = '  -3.39%  '

Set-TextCell "D25" "B25" '84.87'
$ws.Range("E25").Value = '  -1.02%  '

$ws.Range("E26").Value = '  +0.16%  '

Set-TextCell "D27" "B27" '3.00'
$ws.Range("E27").Value = '  -0.26%  '

Set-TextCell "D28" "B28" '9.11'
$ws.Range("E28").Value = '  -2.29%  '

Set-TextCell "D29" "B29" '2.36'
$ws.Range("E29").Value = '  +1.26%  '

Set-TextCell "D30" "B30" '0.121'
$ws.Range("E30").Value = '  +34.27%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell "D31" "B31" '7.00'
$ws.Range("E31").Value = '  -3.17%  '

$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell "D32" "B32" '2.93'
$ws.Range("E32").Value = '  +0.68%  '

Set-TextCell "D33" "B33" '28.16'
$ws.Range("E33").Value = '  -0.60%  '

$ws.Range("E34").Value = '  +0.20%  '

$ws.Range("E35").Value = '  -3.99%  '

Set-TextCell "D36" "B36" '6.47'
$ws.Range("E36").Value = '  -1.16%  '

Set-TextCell "D37" "B37" '55.70'
$ws.Range("E37").Value = '  +1.56%  '

Set-TextCell "D38" "B38" '500.82'
$ws.Range("E38").Value = '  -1.85%  '

Set-TextCell "D39" "B39" '0.0₃0769'
$ws.Range("E39").Value = '  +13.80%  '

Set-TextCell "D40" "B40" '3.09'
$ws.Range("E40").Value = '  +6.84%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell "D41" "B41" '0.132'
$ws.Range("E41").Value = '  +6.49%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D42" "B42" '0.0420'
$ws.Range("E42").Value = '  -1.20%  '

Set-TextCell "D43" "B43" '8.70'
$ws.Range("E43").Value = '  -2.33%  '

Set-TextCell "D44" "B44" '0.297'
$ws.Range("E44").Value = '  -1.97%  '

Set-TextCell "D45" "B45" '2.45'
$ws.Range("E45").Value = '  -0.62%  '

Set-TextCell "D46" "B46" '2.902.28'
$ws.Range("E46").Value = '  -0.77%  '

Set-TextCell "D47" "B47" '28.11'
$ws.Range("E47").Value = '  -2.59%  '

$ws.Range("E48").Value = '  +3.34%  '

$ws.Range("E50").Value = '  -1.21%  '

Set-TextCell "D51" "B51" '122.35'
$ws.Range("E51").Value = '  -0.30%  '
